$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description of "The Earthquake" spell (row 5, column E) to add
# the extra sentence about the grid needing to be empty.
$ws.Range("E5").Value = "Büyü kullanıldığında seçili olan grid tamamen yok olur. Grid üstünde herhangi bir şey olmamalı."

# Move the active selection to E5 (matches the saved sheetView selection).
$ws.Range("E5").Select()
